# Add "Wins" / "Losses" / "Ties" team-record columns (AC, AD, AE) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AB1, style index "1":
# bold, centered, thin border) onto the three new header cells so they match
# the look of the rest of row 1 exactly, then set their text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every team in this sheet had the same 81-81-0 record, so fill it down for
# all 39 data rows (rows 2 through 40).
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 29).Value = 81
    $ws.Cells.Item($r, 30).Value = 81
    $ws.Cells.Item($r, 31).Value = 0
}

Write-Host "done"
